$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.185.21'
$ws.Range("E2").Value = '  +2.62%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.203.48'
$ws.Range("E3").Value = '  +1.20%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.14'
$ws.Range("E5").Value = '  +1.74%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.09'
$ws.Range("E6").Value = '  +4.43%  '

# Row 7
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").Value = '  -1.76%  '

# Row 9
$ws.Range("E9").Value = '  +0.53%  '

# Row 10
$ws.Range("E10").Value = '  +1.02%  '

# Row 11
$ws.Range("E11").Value = '  -1.18%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.755.65'
$ws.Range("E12").Value = '  +1.23%  '

# Row 13
$ws.Range("E13").Value = '  -2.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.86'
$ws.Range("E14").Value = '  +0.43%  '

# Row 15
$ws.Range("E15").Value = '  +0.98%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.186.60'
$ws.Range("E16").Value = '  +2.56%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.212.67'
$ws.Range("E17").Value = '  +1.62%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.29'
$ws.Range("E18").Value = '  +0.43%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.28'
$ws.Range("E19").Value = '  +2.28%  '

# Row 20
$ws.Range("E20").Value = '  +1.22%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '370.99'

# Row 22
$ws.Range("E22").Value = '  -0.04%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.77'
$ws.Range("E24").Value = '  +0.08%  '

# Row 25
$ws.Range("E25").Value = '  +1.73%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.65'
$ws.Range("E26").Value = '  +4.85%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.26%  '

# Row 28
$ws.Range("E28").Value = '  +1.28%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.47'
$ws.Range("E29").Value = '  +0.48%  '

# Row 30
$ws.Range("E30").Value = '  +0.52%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  +1.41%  '

# Row 32
$ws.Range("E32").Value = '  +2.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.57'
$ws.Range("E33").Value = '  +4.23%  '

# Row 34
$ws.Range("E34").Value = '  +2.92%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.14'
$ws.Range("E35").Value = '  +1.39%  '

# Row 36
$ws.Range("E36").Value = '  +2.77%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.46'
$ws.Range("E37").Value = '  +5.36%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.794.33'
$ws.Range("E38").Value = '  +4.27%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0709'
$ws.Range("E39").Value = '  +1.96%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0313'
$ws.Range("E40").Value = '  +8.35%  '

# Row 41
$ws.Range("E41").Value = '  +0.66%  '

# Row 42
$ws.Range("E42").Value = '  -1.51%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.92'
$ws.Range("E43").Value = '  +1.95%  '

# Row 44
$ws.Range("E44").Value = '  -0.50%  '

# Row 45
$ws.Range("E45").Value = '  +0.93%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.242.40'
$ws.Range("E46").Value = '  +1.12%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.985'

# Row 48
$ws.Range("E48").Value = '  -0.95%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.69'
$ws.Range("E49").Value = '  +2.80%  '

# Row 50
$ws.Range("E50").Value = '  +5.88%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.01%  '
